$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 120, shifting existing rows 120-225 down to 121-226.
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row 120 with the new record's data.
$ws.Cells.Item(120, 1).Value = 5
$ws.Cells.Item(120, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(120, 3).Value = "Maule"
$ws.Cells.Item(120, 4).Value = 44658
$ws.Cells.Item(120, 5).Value = 7
$ws.Cells.Item(120, 6).Value = 100112009
$ws.Cells.Item(120, 7).Value = "Acelga"
$ws.Cells.Item(120, 8).Value = "Sin especificar"
$ws.Cells.Item(120, 9).Value = "Primera"
$ws.Cells.Item(120, 10).Value = 400
$ws.Cells.Item(120, 11).Value = 3500
$ws.Cells.Item(120, 12).Value = 3500
$ws.Cells.Item(120, 13).Value = 3500
$ws.Cells.Item(120, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(120, 15).Value = "Región del Maule"
$ws.Cells.Item(120, 16).Value = 875
$ws.Cells.Item(120, 17).Value = 4
$ws.Cells.Item(120, 18).Value = "Hortaliza"

# Match the date cell style used by the rest of column D.
$ws.Cells.Item(120, 4).NumberFormat = $ws.Cells.Item(121, 4).NumberFormat
